$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Clear the value of B4 (Revenue row) while preserving its formatting/style
$ws.Range("B4").ClearContents()
